$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-07 Saturday" "2026-02-08 Sunday"

Replace-Text "70×66=4620" "57×89=5073"
Replace-Text "58×28=1624" "48×25=1200"
Replace-Text "32×23=736" "96×61=5856"
Replace-Text "22×45=990" "32×68=2176"
Replace-Text "76×22=1672" "40×11=440"
Replace-Text "74×96=7104" "20×47=940"
Replace-Text "40×77=3080" "21×66=1386"
Replace-Text "92×55=5060" "83×70=5810"
Replace-Text "49×22=1078" "85×99=8415"
Replace-Text "90×69=6210" "97×66=6402"
Replace-Text "39×45=1755" "51×78=3978"
Replace-Text "65×93=6045" "30×12=360"
Replace-Text "60×31=1860" "36×48=1728"
Replace-Text "56×37=2072" "56×80=4480"
Replace-Text "54×44=2376" "31×19=589"
Replace-Text "45×48=2160" "44×57=2508"
Replace-Text "12×51=612" "81×84=6804"
Replace-Text "60×23=1380" "59×86=5074"
Replace-Text "37×59=2183" "48×56=2688"
Replace-Text "15×55=825" "47×62=2914"
Replace-Text "92×52=4784" "80×47=3760"
Replace-Text "46×94=4324" "98×81=7938"
Replace-Text "57×23=1311" "23×20=460"
Replace-Text "42×34=1428" "57×27=1539"
Replace-Text "90×17=1530" "87×73=6351"
